$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each cryptocurrency row.
# D-column values are written with a leading apostrophe so Excel keeps them
# as text (matching the original inline-string cell type) instead of coercing
# number-looking strings (e.g. "331.07") into floating point numbers. The
# style is then reset to "Normal" so no stray quote-prefix style gets attached.

$ws.Range("D2").Value = "'27.449.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "'1.831.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.86%  "
$ws.Range("D5").Value = "'331.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D7").Value = "'0.4605"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.23%  "
$ws.Range("D8").Value = "'0.3823"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.04%  "
$ws.Range("D9").Value = "'46.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'0.07908"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").Value = "'0.9705"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.51%  "
$ws.Range("D12").Value = "'21.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.68%  "
$ws.Range("D13").Value = "'1.825.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("D14").Value = "'5.884"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "'7.037"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "'87.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "'0.06618"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("D19").Value = "'0.00001028"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("D20").Value = "'17.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "'27.449.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "'5.349"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("D25").Value = "'2.305"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").Value = "'2.043.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.19%  "
$ws.Range("D27").Value = "'157.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "'19.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("D29").Value = "'2.064"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "'5.297"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("D31").Value = "'118.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("D32").Value = "'0.9545"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").Value = "'0.09297"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("D34").Value = "'3.584"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").Value = "'5.248"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("D36").Value = "'1.314"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.76%  "
$ws.Range("D37").Value = "'0.05927"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("D38").Value = "'0.02192"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").Value = "'1.159"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").Value = "'8.041"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "'0.5789"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("D42").Value = "'0.1839"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.91%  "
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").Value = "'1.273"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").Value = "'0.5483"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").Value = "'1.868"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("D48").Value = "'0.06642"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("D49").Value = "'110.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("E50").Value = "  -2.83%  "
$ws.Range("E51").Value = "  -1.05%  "
